$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Range("B104").Select() | Out-Null
